$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.532.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "'1.646.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "'1.004"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").Value = "'302.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").Value = "'0.3835"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("D8").Value = "'0.3601"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("D9").Value = "'50.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("D10").Value = "'0.08164"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").Value = "'1.226"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").Value = "'1.005"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "'22.28"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").Value = "'6.439"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "'7.412"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("D16").Value = "'0.00001217"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("D17").Value = "'1.648.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("D18").Value = "'97.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.49%  "
$ws.Range("D19").Value = "'0.07024"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").Value = "'6.757"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.96%  "
$ws.Range("D21").Value = "'17.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("D22").Value = "'1.006"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").Value = "'12.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.66%  "
$ws.Range("D24").Value = "'23.534.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").Value = "'2.484"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.20%  "
$ws.Range("D26").Value = "'3.028"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.56%  "
$ws.Range("D27").Value = "'21.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("D28").Value = "'153.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("D30").Value = "'133.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").Value = "'1.835.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("D32").Value = "'7.030"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.07%  "
$ws.Range("D33").Value = "'2.251"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.89%  "
$ws.Range("D34").Value = "'12.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.57%  "
$ws.Range("D35").Value = "'1.057"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.26%  "
$ws.Range("D36").Value = "'0.02795"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.85%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").Value = "'0.08792"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.77%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2493"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("D39").Value = "'6.059"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("D40").Value = "'0.06956"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.49%  "
$ws.Range("D41").Value = "'12.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.78%  "
$ws.Range("D42").Value = "'0.6969"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "'1.334"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("D44").Value = "'15.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.95%  "
$ws.Range("D45").Value = "'0.6480"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").Value = "'2.287"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").Value = "'3.953"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Value = "'0.07869"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").Value = "'127.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("D51").Value = "'1.175"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.53%  "
